$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# The _GoBack bookmark currently sits at the end of paragraph 6 ("...maggiore
# numero di articoli pubblicati"). Three new exercises about DELETE
# statements are being added right after that paragraph, and the bookmark
# needs to move to the end of the last of those new paragraphs. Remove the
# existing bookmark now; it is recreated as part of the inserted XML below.
# ---------------------------------------------------------------------------
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

$wns = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'

$p6 = $d.Paragraphs(6)
$insertPos = $p6.Range.End - 1
$insertRange = $d.Range($insertPos, $insertPos)

$newItemsXml = @"
<w:p $wns>
  <w:pPr>
    <w:pStyle w:val="Paragrafoelenco"/>
    <w:numPr><w:ilvl w:val="0"/><w:numId w:val="8"/></w:numPr>
  </w:pPr>
  <w:r><w:t xml:space="preserve">Con riferimento al database IFTS_WS scrivi la query che cancella tutti i post di tipo </w:t></w:r>
  <w:proofErr w:type="spellStart"/>
  <w:r><w:t>revision</w:t></w:r>
  <w:proofErr w:type="spellEnd"/>
</w:p>
<w:p $wns>
  <w:pPr>
    <w:pStyle w:val="Paragrafoelenco"/>
    <w:numPr><w:ilvl w:val="0"/><w:numId w:val="8"/></w:numPr>
  </w:pPr>
  <w:r><w:t>Con riferimento al database IFTS scrivi la query che cancella tutti i clienti della città di Parma</w:t></w:r>
</w:p>
<w:p $wns>
  <w:pPr>
    <w:pStyle w:val="Paragrafoelenco"/>
    <w:numPr><w:ilvl w:val="0"/><w:numId w:val="8"/></w:numPr>
  </w:pPr>
  <w:r><w:t>Con riferimento al database IFTS scrivi la query che cancella tutt</w:t></w:r>
  <w:r><w:t>e le prenotazioni che hanno data di arrivo precedente al 01/02/2015</w:t></w:r>
  <w:r><w:t xml:space="preserve"> e importo inferiore a 100€</w:t></w:r>
  <w:bookmarkStart w:id="0" w:name="_GoBack"/>
  <w:bookmarkEnd w:id="0"/>
</w:p>
"@

$insertRange.InsertXML($newItemsXml)

# ---------------------------------------------------------------------------
# A trailing empty paragraph, still styled "Paragrafoelenco" but without list
# numbering, is added right before the document's final blank paragraph.
# ---------------------------------------------------------------------------
$lastPara = $d.Paragraphs($d.Paragraphs.Count)
$lastPos = $lastPara.Range.Start
$lastRange = $d.Range($lastPos, $lastPos)
$emptyItemXml = "<w:p $wns><w:pPr><w:pStyle w:val=`"Paragrafoelenco`"/></w:pPr></w:p>"
$lastRange.InsertXML($emptyItemXml)
